# Actualización automática de tasas-transfi.xlsx
#
# 1) Hoja1!A1 - update the two "Binance" conversion bullet lines inside the
#    daily-update text block (keeps every other line, including blank
#    lines, untouched).
# 2) tasas!N10 / tasas!O10 / tasas!N12 / tasas!O12 - refreshed rate figures.

$wb = $excel.ActiveWorkbook

# --- 1) Update the "Conversión del día" note on Hoja1 -----------------
$notes = $wb.Worksheets.Item("Hoja1")
$cell = $notes.Range("A1")
$text = $cell.Text
$text = $text.Replace("1000 Bs = 3.17 = 12174.33 pesos", "1000 Bs = 3.22 = 12352.8 pesos")
$text = $text.Replace("12174.33 pesos = 3.16 = 958.63 Bs", "12352.8 pesos = 3.2 = 965.53 Bs")
$cell.Value = $text

# --- 2) Refresh the rate figures on the "tasas" sheet ------------------
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 311
$tasas.Range("O10").Value = 3841.72
$tasas.Range("N12").Value = 3864.99
$tasas.Range("O12").Value = 302.1
